# "Debugged default options for 'General Settings'":
# Sheet2 holds the default/general-settings row (row 2). A stray debug
# value of 41 had landed in Stop Row (C2), 3 in Skip Row (D2), and the
# Transpose flag (M2) was incorrectly defaulted to "Yes" - clear all
# three so the row goes back to its intended blank defaults.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # Sheet2 is tabSelected, so it's already active

$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("M2").ClearContents()

# Move / save the selection to C2, matching the saved view state.
$ws.Range("C2").Select()
